$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0003714022599530242
$ws.Range("C2").Value = 0.004309184025731883
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 250.0725307070909
